# Rename the embedded logo pictures' display names:
#   - the two Pearson Edexcel logo pictures (in the default & first-page
#     footers) go from "image1.png" to "image2.png"
#   - the BTEC logo picture (in the first-page header) goes from
#     "image2.jpg" to "image1.jpg"
#
# InlineShape objects don't expose a writable .Name in the Word object
# model, so each picture is temporarily converted to a floating Shape
# (which does expose .Name), renamed, then converted back to an inline
# picture so the layout/anchoring is unchanged.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($range, $newName) {
    if ($range.InlineShapes.Count -gt 0) {
        $inlineShape = $range.InlineShapes(1)
        $floating = $inlineShape.ConvertToShape()
        $floating.Name = $newName
        [void]$floating.ConvertToInlineShape()
    }
}

# Pearson Edexcel logo -> image2.png (both footers that carry it)
Rename-InlinePicture $sec.Footers(1).Range "image2.png"
Rename-InlinePicture $sec.Footers(2).Range "image2.png"

# BTEC logo -> image1.jpg (first-page header)
Rename-InlinePicture $sec.Headers(2).Range "image1.jpg"
